# edit.ps1
# Implements: insert a new "CAN VIEW" / can_view permission row into the
# INDO_CMS_MENU_PERMISSION block of sheet "INDO_CMS_TEMPLATE_DETAIL",
# right after "ROLE ID" and before "CAN INSERT". This pushes all rows
# below it down by one (the template header / template detail blocks
# simply shift down one row, keeping their own internal sequence numbers).
# The sequence numbers (column G) for the permission rows from ROLE ID
# onward are bumped by +1 to account for the newly inserted row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("INDO_CMS_TEMPLATE_DETAIL")
$ws2 = $wb.Worksheets.Item("INDO_CMS_TEMPLATE_DETAIL_2")

# --- Insert a new row at row 29 (pushes CAN INSERT..CAN IMPORT, the
# TEMPLATE_HEADER block and the TEMPLATE_DETAIL block all down by one row) ---
$ws1.Rows.Item(29).Insert()

# --- Fill in the newly inserted row 29 with the CAN VIEW permission ---
$ws1.Range("B29").Value = "INDO_CMS_MENU_PERMISSION"
$ws1.Range("C29").Value = "CAN VIEW"
$ws1.Range("D29").Value = "can_view"
$ws1.Range("E29").Value = "can_view"
$ws1.Range("F29").Value = "STRING"
$ws1.Range("G29").Value = 4
$ws1.Range("H29").Value = 0
$ws1.Range("I29").Value = 1
$ws1.Range("J29").Value = 1
$ws1.Range("K29").Value = 1
$ws1.Range("M29").Formula = '="(''"&B29&"'',''"&C29&"'',''"&D29&"'',''"&E29&"'',''"&F29&"'',''"&G29&"'',''"&H29&"'',''"&I29&"'',''"&J29&"'',''"&K29&"'');"'
$ws1.Range("N29").Formula = '="INSERT INTO INDO_CMS_TEMPLATE_DETAIL (template_code,web_column,database_column,query_column,data_type,sequence,is_primary,is_show,is_create,is_edit) VALUES"'
$ws1.Range("P29").Formula = "=N29&M29"

# --- Bump the sequence (column G) for ROLE ID (row 28) and for the rows
# that shifted down (CAN INSERT..CAN IMPORT, now rows 30..34) ---
$ws1.Range("G28").Value = 3
$ws1.Range("G30").Value = 5
$ws1.Range("G31").Value = 6
$ws1.Range("G32").Value = 7
$ws1.Range("G33").Value = 8
$ws1.Range("G34").Value = 9

# --- View state: selection on the active sheet moves to P26:P34 ---
$ws1.Range("P26:P34").Select()

# --- Best-effort: mirror the lingering multi-area selection state that
# was left on the second sheet (sqref "P26:P34 A1") ---
$ws2.Range("A1").Select()
